$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = 45411
$ws.Cells.Item(2, 2).Value = 7.8125
$ws.Cells.Item(2, 3).Value = 3.5
$ws.Cells.Item(2, 4).Value = 10.5
$ws.Cells.Item(2, 5).Value = 88.175
$ws.Cells.Item(2, 6).Value = 1011.051944743122
$ws.Cells.Item(2, 7).Value = 0
$ws.Cells.Item(2, 8).Value = 0
$ws.Cells.Item(2, 9).Value = 4
$ws.Cells.Item(2, 10).Value = 29
$ws.Cells.Item(2, 11).Value = 5.791377543624169
$ws.Cells.Item(2, 12).Value = 5.677516843120022
$ws.Cells.Item(2, 13).Value = 5.595194113481074
$ws.Cells.Item(2, 14).Value = 5.151161093946595
$ws.Cells.Item(2, 15).Value = 4.521395945090218
$ws.Cells.Item(2, 16).Value = 3.420377619647766

# Row 3
$ws.Cells.Item(3, 1).Value = 45412
$ws.Cells.Item(3, 2).Value = 9.725
$ws.Cells.Item(3, 3).Value = 3.5
$ws.Cells.Item(3, 4).Value = 14.3
$ws.Cells.Item(3, 5).Value = 89.21249999999999
$ws.Cells.Item(3, 6).Value = 1015.850915142929
$ws.Cells.Item(3, 7).Value = 0
$ws.Cells.Item(3, 8).Value = 0
$ws.Cells.Item(3, 9).Value = 4
$ws.Cells.Item(3, 10).Value = 30
$ws.Cells.Item(3, 11).Value = 7.353811737318463
$ws.Cells.Item(3, 12).Value = 7.25727072305879
$ws.Cells.Item(3, 13).Value = 7.142052262331761
$ws.Cells.Item(3, 14).Value = 6.515144020162255
$ws.Cells.Item(3, 15).Value = 5.39911823591947
$ws.Cells.Item(3, 16).Value = 4.15617515858971

# Row 4
$ws.Cells.Item(4, 1).Value = 45413
$ws.Cells.Item(4, 2).Value = 11.4875
$ws.Cells.Item(4, 3).Value = 6.2
$ws.Cells.Item(4, 4).Value = 16.3
$ws.Cells.Item(4, 5).Value = 68.6875
$ws.Cells.Item(4, 6).Value = 1014.289883072826
$ws.Cells.Item(4, 7).Value = 0
$ws.Cells.Item(4, 8).Value = 0
$ws.Cells.Item(4, 9).Value = 5
$ws.Cells.Item(4, 10).Value = 1
$ws.Cells.Item(4, 11).Value = 9.746383534853404
$ws.Cells.Item(4, 12).Value = 9.583914022080549
$ws.Cells.Item(4, 13).Value = 9.455973997361156
$ws.Cells.Item(4, 14).Value = 9.000052572493837
$ws.Cells.Item(4, 15).Value = 7.964533447735224
$ws.Cells.Item(4, 16).Value = 6.161662057835041

# Row 5
$ws.Cells.Item(5, 1).Value = 45414
$ws.Cells.Item(5, 2).Value = 14.225
$ws.Cells.Item(5, 3).Value = 8.2
$ws.Cells.Item(5, 4).Value = 19.8
$ws.Cells.Item(5, 5).Value = 67.6
$ws.Cells.Item(5, 6).Value = 1010.095566015488
$ws.Cells.Item(5, 7).Value = 0
$ws.Cells.Item(5, 8).Value = 0
$ws.Cells.Item(5, 9).Value = 5
$ws.Cells.Item(5, 10).Value = 2
$ws.Cells.Item(5, 11).Value = 11.17226086676881
$ws.Cells.Item(5, 12).Value = 10.96551109576654
$ws.Cells.Item(5, 13).Value = 10.79630447786615
$ws.Cells.Item(5, 14).Value = 10.14202901156393
$ws.Cells.Item(5, 15).Value = 9.009832640559624
$ws.Cells.Item(5, 16).Value = 6.941975406800601

# Row 6
$ws.Cells.Item(6, 1).Value = 45415
$ws.Cells.Item(6, 2).Value = 14.4
$ws.Cells.Item(6, 3).Value = 7.2
$ws.Cells.Item(6, 4).Value = 19.7
$ws.Cells.Item(6, 5).Value = 69
$ws.Cells.Item(6, 6).Value = 1006.789173293079
$ws.Cells.Item(6, 7).Value = 0
$ws.Cells.Item(6, 8).Value = 0
$ws.Cells.Item(6, 9).Value = 5
$ws.Cells.Item(6, 10).Value = 3
$ws.Cells.Item(6, 11).Value = 11.20551190392998
$ws.Cells.Item(6, 12).Value = 11.02469714287919
$ws.Cells.Item(6, 13).Value = 10.85969954574997
$ws.Cells.Item(6, 14).Value = 10.17898013612291
$ws.Cells.Item(6, 15).Value = 9.092596766998762
$ws.Cells.Item(6, 16).Value = 6.889598149941348

# Row 7
$ws.Cells.Item(7, 1).Value = 45416
$ws.Cells.Item(7, 2).Value = 14.475
$ws.Cells.Item(7, 3).Value = 8.3
$ws.Cells.Item(7, 4).Value = 19.5
$ws.Cells.Item(7, 5).Value = 64.475
$ws.Cells.Item(7, 6).Value = 1005.061298545689
$ws.Cells.Item(7, 7).Value = 0
$ws.Cells.Item(7, 8).Value = 0
$ws.Cells.Item(7, 9).Value = 5
$ws.Cells.Item(7, 10).Value = 4
$ws.Cells.Item(7, 11).Value = 11.13823833034559
$ws.Cells.Item(7, 12).Value = 10.94169111688645
$ws.Cells.Item(7, 13).Value = 10.76945745094702
$ws.Cells.Item(7, 14).Value = 10.13997430765606
$ws.Cells.Item(7, 15).Value = 9.009780719034444
$ws.Cells.Item(7, 16).Value = 6.94197507760016

# Row 8
$ws.Cells.Item(8, 1).Value = 45417
$ws.Cells.Item(8, 2).Value = 12.625
$ws.Cells.Item(8, 3).Value = 7
$ws.Cells.Item(8, 4).Value = 16.9
$ws.Cells.Item(8, 5).Value = 66.5
$ws.Cells.Item(8, 6).Value = 1008.178773733112
$ws.Cells.Item(8, 7).Value = 0
$ws.Cells.Item(8, 8).Value = 0
$ws.Cells.Item(8, 9).Value = 5
$ws.Cells.Item(8, 10).Value = 5
$ws.Cells.Item(8, 11).Value = 10.71249358589029
$ws.Cells.Item(8, 12).Value = 10.56510486885183
$ws.Cells.Item(8, 13).Value = 10.57740809076677
$ws.Cells.Item(8, 14).Value = 10.05278383724556
$ws.Cells.Item(8, 15).Value = 9.007577451601327
$ws.Cells.Item(8, 16).Value = 6.941961108122077

# Row 9
$ws.Cells.Item(9, 1).Value = 45418
$ws.Cells.Item(9, 2).Value = 11.4
$ws.Cells.Item(9, 3).Value = 5.8
$ws.Cells.Item(9, 4).Value = 15.7
$ws.Cells.Item(9, 5).Value = 70.65
$ws.Cells.Item(9, 6).Value = 1009.267751964902
$ws.Cells.Item(9, 7).Value = 0
$ws.Cells.Item(9, 8).Value = 0
$ws.Cells.Item(9, 9).Value = 5
$ws.Cells.Item(9, 10).Value = 6
$ws.Cells.Item(9, 11).Value = 9.878128359683625
$ws.Cells.Item(9, 12).Value = 9.802533281050177
$ws.Cells.Item(9, 13).Value = 9.723369469326094
$ws.Cells.Item(9, 14).Value = 9.087083711455918
$ws.Cells.Item(9, 15).Value = 8.157014119254471
$ws.Cells.Item(9, 16).Value = 5.909927118985828

# Row 10
$ws.Cells.Item(10, 1).Value = 45419
$ws.Cells.Item(10, 2).Value = 11.225
$ws.Cells.Item(10, 3).Value = 5.6
$ws.Cells.Item(10, 4).Value = 15
$ws.Cells.Item(10, 5).Value = 66.75
$ws.Cells.Item(10, 6).Value = 1010.941929630209
$ws.Cells.Item(10, 7).Value = 0
$ws.Cells.Item(10, 8).Value = 0
$ws.Cells.Item(10, 9).Value = 5
$ws.Cells.Item(10, 10).Value = 7
$ws.Cells.Item(10, 11).Value = 9.742964699829432
$ws.Cells.Item(10, 12).Value = 9.583144016141574
$ws.Cells.Item(10, 13).Value = 9.455724787595667
$ws.Cells.Item(10, 14).Value = 9.00003349953352
$ws.Cells.Item(10, 15).Value = 7.964532965769338
$ws.Cells.Item(10, 16).Value = 6.16166205477921

# Row 11
$ws.Cells.Item(11, 1).Value = 45420
$ws.Cells.Item(11, 2).Value = 11.075
$ws.Cells.Item(11, 3).Value = 6.8
$ws.Cells.Item(11, 4).Value = 13.9
$ws.Cells.Item(11, 5).Value = 72
$ws.Cells.Item(11, 6).Value = 1011.306800285849
$ws.Cells.Item(11, 7).Value = 0
$ws.Cells.Item(11, 8).Value = 0
$ws.Cells.Item(11, 9).Value = 5
$ws.Cells.Item(11, 10).Value = 8
$ws.Cells.Item(11, 11).Value = 10.03086049699492
$ws.Cells.Item(11, 12).Value = 9.901546065648002
$ws.Cells.Item(11, 13).Value = 9.823226181482399
$ws.Cells.Item(11, 14).Value = 9.454338332837597
$ws.Cells.Item(11, 15).Value = 8.305700911365102
$ws.Cells.Item(11, 16).Value = 6.353378191254905

# Row 12
$ws.Cells.Item(12, 1).Value = 45421
$ws.Cells.Item(12, 2).Value = 11.83333333333333
$ws.Cells.Item(12, 3).Value = 7
$ws.Cells.Item(12, 4).Value = 16.4
$ws.Cells.Item(12, 5).Value = 73.96666666666667
$ws.Cells.Item(12, 6).Value = 1010.90011243065
$ws.Cells.Item(12, 7).Value = 0
$ws.Cells.Item(12, 8).Value = 0
$ws.Cells.Item(12, 9).Value = 5
$ws.Cells.Item(12, 10).Value = 9
$ws.Cells.Item(12, 11).Value = 10.92125957150249
$ws.Cells.Item(12, 12).Value = 10.75432625041842
$ws.Cells.Item(12, 13).Value = 10.65353229403615
$ws.Cells.Item(12, 14).Value = 10.08450872932704
$ws.Cells.Item(12, 15).Value = 9.008379126614813
$ws.Cells.Item(12, 16).Value = 6.941966191019125

